# Generate Report for Handoff
# Replace the old handoff UUID (a1cdd7fa-8980-4763-be32-a7e7acd7bb17) with the
# new one (1f05dc6c-ec38-449a-89f1-5d2a15e41168), refresh the handoff/xliff
# timestamps and xliff hash, and update the hyperlink display text on every
# sheet to match (the hyperlink target addresses themselves stay the same).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # Overview
$ws2 = $wb.Worksheets.Item(2)   # zh-cn
$ws3 = $wb.Worksheets.Item(3)   # de-de

$oldGuid = "a1cdd7fa-8980-4763-be32-a7e7acd7bb17"
$newGuid = "1f05dc6c-ec38-449a-89f1-5d2a15e41168"

$oldHash = "cb7940d374ff5c2632c4bcbcb167ed379c22b324"
$newHash = "c7ed9542314e0f2fa710694eae48b7f3013bef12"

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$ws1.Range("A2").Value = $newGuid + ".md"
$ws1.Range("B2").Value = "e2e\" + $newGuid + ".md"
$ws1.Range("G2").Value = "2016-09-06 15:14:41"

$ws1.Range("B2").Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8d43581d1f9319f670a699bdf503316172dc1810/e2e/" + $oldGuid + ".md", "", "", "e2e\" + $newGuid + ".md")

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$ws2.Range("A2").Value = $newGuid + ".md"
$ws2.Range("G2").Value = $newGuid + "." + $newHash + ".zh-cn.xlf"
$ws2.Range("H2").Value = "2016-09-06 15:14:27"

$ws2.Range("A2").Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8d43581d1f9319f670a699bdf503316172dc1810/e2e/" + $oldGuid + ".md", "", "", $newGuid + ".md")

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$ws3.Range("A2").Value = $newGuid + ".md"
$ws3.Range("G2").Value = $newGuid + "." + $newHash + ".de-de.xlf"
$ws3.Range("H2").Value = "2016-09-06 15:14:41"

$ws3.Range("A2").Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8d43581d1f9319f670a699bdf503316172dc1810/e2e/" + $oldGuid + ".md", "", "", $newGuid + ".md")
